# Natmi following Dr Hou advice
# Updates the LR-pairs sheet (Il10-Il10rb) with a new sending cluster ("sCs")
# and refreshed statistics for existing ECs rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il10"
$ws.Cells.Item(2, 3).Value = "Il10rb"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 37.63846333333333
$ws.Cells.Item(2, 8).Value = 112.91539
$ws.Cells.Item(2, 9).Value = 0.9976844064794065
$ws.Cells.Item(2, 10).Value = 0.9976844064794066
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 81.038515
$ws.Cells.Item(2, 14).Value = 243.115545
$ws.Cells.Item(2, 15).Value = 0.8313546010859307
$ws.Cells.Item(2, 16).Value = 0.8313546010859306
$ws.Cells.Item(2, 17).Value = 3050.165175415284
$ws.Cells.Item(2, 18).Value = 27451.48657873755
$ws.Cells.Item(2, 19).Value = 0.8294295217583405
$ws.Cells.Item(2, 20).Value = 0.8294295217583405

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il10"
$ws.Cells.Item(3, 3).Value = "Il10rb"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 37.63846333333333
$ws.Cells.Item(3, 8).Value = 112.91539
$ws.Cells.Item(3, 9).Value = 0.9976844064794065
$ws.Cells.Item(3, 10).Value = 0.9976844064794066
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 12.88200833333333
$ws.Cells.Item(3, 14).Value = 38.646025
$ws.Cells.Item(3, 15).Value = 0.1321534198787326
$ws.Cells.Item(3, 16).Value = 0.1321534198787326
$ws.Cells.Item(3, 17).Value = 484.8589983138611
$ws.Cells.Item(3, 18).Value = 4363.730984824751
$ws.Cells.Item(3, 19).Value = 0.1318474062759371
$ws.Cells.Item(3, 20).Value = 0.1318474062759371

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il10"
$ws.Cells.Item(4, 3).Value = "Il10rb"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 37.63846333333333
$ws.Cells.Item(4, 8).Value = 112.91539
$ws.Cells.Item(4, 9).Value = 0.9976844064794065
$ws.Cells.Item(4, 10).Value = 0.9976844064794066
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.557153333333333
$ws.Cells.Item(4, 14).Value = 10.67146
$ws.Cells.Item(4, 15).Value = 0.03649197903533674
$ws.Cells.Item(4, 16).Value = 0.03649197903533673
$ws.Cells.Item(4, 17).Value = 133.8857853077111
$ws.Cells.Item(4, 18).Value = 1204.9720677694
$ws.Cells.Item(4, 19).Value = 0.03640747844512888
$ws.Cells.Item(4, 20).Value = 0.03640747844512888

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il10"
$ws.Cells.Item(5, 3).Value = "Il10rb"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.08735766666666667
$ws.Cells.Item(5, 8).Value = 0.262073
$ws.Cells.Item(5, 9).Value = 0.002315593520593406
$ws.Cells.Item(5, 10).Value = 0.002315593520593406
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 81.038515
$ws.Cells.Item(5, 14).Value = 243.115545
$ws.Cells.Item(5, 15).Value = 0.8313546010859307
$ws.Cells.Item(5, 16).Value = 0.8313546010859306
$ws.Cells.Item(5, 17).Value = 7.079335580531667
$ws.Cells.Item(5, 18).Value = 63.714020224785
$ws.Cells.Item(5, 19).Value = 0.001925079327590097
$ws.Cells.Item(5, 20).Value = 0.001925079327590097

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il10"
$ws.Cells.Item(6, 3).Value = "Il10rb"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.08735766666666667
$ws.Cells.Item(6, 8).Value = 0.262073
$ws.Cells.Item(6, 9).Value = 0.002315593520593406
$ws.Cells.Item(6, 10).Value = 0.002315593520593406
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 12.88200833333333
$ws.Cells.Item(6, 14).Value = 38.646025
$ws.Cells.Item(6, 15).Value = 0.1321534198787326
$ws.Cells.Item(6, 16).Value = 0.1321534198787326
$ws.Cells.Item(6, 17).Value = 1.125342189980556
$ws.Cells.Item(6, 18).Value = 10.128079709825
$ws.Cells.Item(6, 19).Value = 0.000306013602795453
$ws.Cells.Item(6, 20).Value = 0.000306013602795453

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Il10"
$ws.Cells.Item(7, 3).Value = "Il10rb"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.08735766666666667
$ws.Cells.Item(7, 8).Value = 0.262073
$ws.Cells.Item(7, 9).Value = 0.002315593520593406
$ws.Cells.Item(7, 10).Value = 0.002315593520593406
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.557153333333333
$ws.Cells.Item(7, 14).Value = 10.67146
$ws.Cells.Item(7, 15).Value = 0.03649197903533674
$ws.Cells.Item(7, 16).Value = 0.03649197903533673
$ws.Cells.Item(7, 17).Value = 0.3107446151755556
$ws.Cells.Item(7, 18).Value = 2.79670153658
$ws.Cells.Item(7, 19).Value = 0.00008450059020785619
$ws.Cells.Item(7, 20).Value = 0.00008450059020785616

